# Add the new "Sheet1" worksheet after the existing "Blad1" sheet.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# Enter the labels in the same order the original author typed them so the
# shared-string table comes out in the same sequence.
$ws.Range("C3").Value = "Vsource+"
$ws.Range("C4").Value = "Vsource-"
$ws.Range("C5").Value = "Bit setting (0 - 32)"
$ws.Range("F3").Value = "Rshunt"
$ws.Range("F5").Value = "Current"
$ws.Range("F7").Value = "Comp+ input"
$ws.Range("C7").Value = "Comp- input"
$ws.Range("C6").Value = "DAC Vout"
$ws.Range("C11").Value = "Vsource+ can be externall ref (2.5V) or Vdd(5.0V) or Fvr(4.096V)"
$ws.Range("F4").Value = "Gain INA326"

# Numeric inputs / formulas, column D.
$ws.Range("D3").Value = 4.096
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 22
$ws.Range("D6").Formula = "=((D3-D4)*(D5/32))+ D4"
$ws.Range("D7").Formula = "=D6"

# Numeric inputs / formulas, column G.
$ws.Range("G3").Value = 0.025
$ws.Range("G4").Value = 5000
$ws.Range("G5").Value = 0.003
$ws.Range("G7").Formula = "=G3*G5*G4+2.5"

# Number formatting - three decimal places on all the calc cells.
$ws.Range("D3:D7").NumberFormat = "0.000"
$ws.Range("G3").NumberFormat = "0.000"
$ws.Range("G7").NumberFormat = "0.000"

# "Good" (green) cell style highlighting the user-editable inputs.
$ws.Range("C5").Style = "Good"
$ws.Range("G4").Style = "Good"
$ws.Range("G4").NumberFormat = "0.000"
$ws.Range("G5").Style = "Good"
$ws.Range("G5").NumberFormat = "0.000"

# Column widths observed in the authored sheet (best-fit on the label
# columns). The stored OOXML width is `ColumnWidth + ~0.833`, so back the
# COM value off to land on the authored width as closely as possible.
$ws.Columns.Item(2).ColumnWidth = 11.1666666667
$ws.Columns.Item(3).ColumnWidth = 15.7369791667
$ws.Columns.Item(6).ColumnWidth = 11.5924479167

# Select cell I10 on the new sheet, and make it the active sheet/tab.
$null = $ws.Range("I10").Select()
$null = $ws.Activate()

Write-Output "done"
